$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.922.13'
$ws.Range('E2').Value = '  +6.87%  '
$ws.Range('D3').Value = '2.311.33'
$ws.Range('E3').Value = '  +5.42%  '
$ws.Range('E4').Value = '  -0.59%  '
$range = $ws.Range('D5')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '298.77'
$range.Style = $origStyle
$ws.Range('E5').Value = '  +1.77%  '
$range = $ws.Range('D6')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '99.05'
$range.Style = $origStyle
$ws.Range('E6').Value = '  +12.57%  '
$range = $ws.Range('D7')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '0.573'
$range.Style = $origStyle
$ws.Range('E7').Value = '  +1.59%  '
$range = $ws.Range('D8')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '0.999'
$range.Style = $origStyle
$ws.Range('E8').Value = '  -0.47%  '
$ws.Range('E9').Value = '  +10.57%  '
$range = $ws.Range('D10')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '35.72'
$range.Style = $origStyle
$ws.Range('E10').Value = '  +11.31%  '
$ws.Range('E11').Value = '  +4.74%  '
$ws.Range('E12').Value = '  +8.98%  '
$ws.Range('E13').Value = '  +1.93%  '
$ws.Range('D14').Value = '2.662.94'
$ws.Range('E14').Value = '  +5.43%  '
$ws.Range('D15').Value = '2.313.11'
$ws.Range('E15').Value = '  +3.07%  '
$range = $ws.Range('D16')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '13.99'
$range.Style = $origStyle
$ws.Range('E16').Value = '  +8.61%  '
$range = $ws.Range('D17')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '0.816'
$range.Style = $origStyle
$ws.Range('E17').Value = '  +7.16%  '
$ws.Range('D18').Value = '46.804.64'
$ws.Range('E18').Value = '  +7.76%  '
$range = $ws.Range('D19')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '13.19'
$range.Style = $origStyle
$ws.Range('E19').Value = '  +23.73%  '
$ws.Range('D20').Value = '0.0₃0941'
$ws.Range('E20').Value = '  +7.40%  '
$ws.Range('E21').Value = '  +6.08%  '
$range = $ws.Range('D22')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '66.83'
$range.Style = $origStyle
$ws.Range('E22').Value = '  +6.91%  '
$range = $ws.Range('D23')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '248.14'
$range.Style = $origStyle
$ws.Range('E23').Value = '  +8.43%  '
$ws.Range('E24').Value = '  +5.98%  '
$range = $ws.Range('D25')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '1.97'
$range.Style = $origStyle
$ws.Range('E25').Value = '  +9.61%  '
$ws.Range('E26').Value = '  -0.24%  '
$range = $ws.Range('D27')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '42.77'
$range.Style = $origStyle
$ws.Range('E27').Value = '  +21.84%  '
$range = $ws.Range('D28')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '2.26'
$range.Style = $origStyle
$ws.Range('E28').Value = '  +1.76%  '
$range = $ws.Range('D29')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '9.89'
$range.Style = $origStyle
$ws.Range('E29').Value = '  +8.19%  '
$range = $ws.Range('D30')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '20.20'
$range.Style = $origStyle
$ws.Range('E30').Value = '  +6.11%  '
$range = $ws.Range('D31')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '5.76'
$range.Style = $origStyle
$ws.Range('E31').Value = '  +9.34%  '
$range = $ws.Range('D32')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '147.12'
$range.Style = $origStyle
$ws.Range('E32').Value = '  +1.19%  '
$ws.Range('E33').Value = '  +10.69%  '
$ws.Range('E34').Value = '  +4.42%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$range = $ws.Range('D35')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '0.113'
$range.Style = $origStyle
$ws.Range('E35').Value = '  +12.18%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$range = $ws.Range('D36')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '3.10'
$range.Style = $origStyle
$ws.Range('E36').Value = '  +8.98%  '
$ws.Range('E37').Value = '  +3.30%  '
$ws.Range('E38').Value = '  +9.79%  '
$range = $ws.Range('D39')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '15.76'
$range.Style = $origStyle
$ws.Range('E39').Value = '  +19.81%  '
$range = $ws.Range('D40')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '4.03'
$range.Style = $origStyle
$ws.Range('E40').Value = '  +16.11%  '
$range = $ws.Range('D41')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '3.41'
$range.Style = $origStyle
$ws.Range('E41').Value = '  +12.51%  '
$ws.Range('E42').Value = '  +10.58%  '
$ws.Range('E43').Value = '  -0.56%  '
$ws.Range('E44').Value = '  +22.48%  '
$ws.Range('D45').Value = '1.841.81'
$ws.Range('E45').Value = '  +5.61%  '
$range = $ws.Range('D46')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '90.54'
$range.Style = $origStyle
$ws.Range('E46').Value = '  +25.20%  '
$ws.Range('E47').Value = '  +17.19%  '
$range = $ws.Range('D48')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '75.96'
$range.Style = $origStyle
$ws.Range('E48').Value = '  +15.56%  '
$ws.Range('E49').Value = '  +10.90%  '
$range = $ws.Range('D50')
$origStyle = $range.Style
$range.NumberFormat = '@'
$range.Value = '97.09'
$range.Style = $origStyle
$ws.Range('E50').Value = '  +6.88%  '
$ws.Range('E51').Value = '  +12.09%  '
